# "addind preparer to sheet"
# Update the "purpose" column (E2:E29) from "S.GISH" to the new value
# "fullRNASEQ", and leave the sheet's selection on the range the author
# ended up with after making the edit (D30:F34, active cell D30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E29").Value = "fullRNASEQ"

$ws.Range("D30:F34").Select()
